# "Mas mediciones de temperatura"
# The raw-data sheet ("Datos crudos") received two additional temperature
# readings (rows 36 and 37) and every existing reading's timestamp/value was
# refreshed to a new measurement run (2023-12-12 instead of 2023-12-08), and
# the "current row" pointer (H2) moved from 11 to 20.  Everything else
# (the COUNT() in H3, the "Datos válidos" sheet, and both charts' cached
# series) is formula/cache driven and recomputes on its own once the raw
# data changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos crudos")
$ws2 = $wb.Worksheets.Item("Datos válidos")

# New timestamp (column C) + temperature (column E) readings for rows 2-37.
$rows = @(
  @{Row=2;  Ts="2023-12-12 04:04:36"; E=26.722027972027899},
  @{Row=3;  Ts="2023-12-12 04:05:35"; E=26.328671328671302},
  @{Row=4;  Ts="2023-12-12 04:06:35"; E=25.935314685314601},
  @{Row=5;  Ts="2023-12-12 04:07:35"; E=25.673076923076898},
  @{Row=6;  Ts="2023-12-12 04:08:35"; E=25.279720279720198},
  @{Row=7;  Ts="2023-12-12 04:09:35"; E=25.148601398601301},
  @{Row=8;  Ts="2023-12-12 04:10:34"; E=25.017482517482499},
  @{Row=9;  Ts="2023-12-12 04:11:34"; E=25.017482517482499},
  @{Row=10; Ts="2023-12-12 04:12:34"; E=24.886363636363601},
  @{Row=11; Ts="2023-12-12 04:13:34"; E=25.017482517482499},
  @{Row=12; Ts="2023-12-12 04:14:34"; E=24.624125874125799},
  @{Row=13; Ts="2023-12-12 04:15:34"; E=24.886363636363601},
  @{Row=14; Ts="2023-12-12 04:16:33"; E=24.493006993006901},
  @{Row=15; Ts="2023-12-12 04:17:33"; E=24.7552447552447},
  @{Row=16; Ts="2023-12-12 04:18:33"; E=24.493006993006901},
  @{Row=17; Ts="2023-12-12 04:19:33"; E=24.624125874125799},
  @{Row=18; Ts="2023-12-12 04:20:33"; E=24.493006993006901},
  @{Row=19; Ts="2023-12-12 04:21:33"; E=24.493006993006901},
  @{Row=20; Ts="2023-12-12 04:22:32"; E=24.361888111888099},
  @{Row=21; Ts="2023-12-12 04:23:32"; E=24.493006993006901},
  @{Row=22; Ts="2023-12-12 04:24:32"; E=24.493006993006901},
  @{Row=23; Ts="2023-12-12 04:25:32"; E=24.361888111888099},
  @{Row=24; Ts="2023-12-12 04:26:32"; E=24.361888111888099},
  @{Row=25; Ts="2023-12-12 04:27:32"; E=24.361888111888099},
  @{Row=26; Ts="2023-12-12 04:28:31"; E=24.493006993006901},
  @{Row=27; Ts="2023-12-12 04:29:31"; E=24.230769230769202},
  @{Row=28; Ts="2023-12-12 04:30:31"; E=24.493006993006901},
  @{Row=29; Ts="2023-12-12 04:31:31"; E=24.361888111888099},
  @{Row=30; Ts="2023-12-12 04:32:31"; E=24.361888111888099},
  @{Row=31; Ts="2023-12-12 04:33:31"; E=24.230769230769202},
  @{Row=32; Ts="2023-12-12 04:34:30"; E=24.361888111888099},
  @{Row=33; Ts="2023-12-12 04:35:30"; E=24.230769230769202},
  @{Row=34; Ts="2023-12-12 04:36:30"; E=24.493006993006901},
  @{Row=35; Ts="2023-12-12 04:37:30"; E=24.624125874125799},
  @{Row=36; Ts="2023-12-12 04:38:30"; E=24.230769230769202},
  @{Row=37; Ts="2023-12-12 04:39:29"; E=24.361888111888099}
)

foreach ($r in $rows) {
  $row = $r.Row
  # Rows 36/37 are brand-new measurement rows -> fill in A/B/D too.
  if ($row -gt 35) {
    $ws.Range("A" + $row).Value = 6
    $ws.Range("B" + $row).Value = 24
    $ws.Range("D" + $row).Value = 0
  }
  $ws.Range("C" + $row).Value = $r.Ts
  $ws.Range("C" + $row).NumberFormat = "@"
  $ws.Range("E" + $row).Value = $r.E
}

# The "current row" pointer used by INDEX()/the "Datos válidos" sheet.
$ws.Range("H2").Value = 20

# Selections recorded in the saved file (cosmetic, but part of the diff).
$ws.Range("O8").Select()
$ws2.Range("I26").Select()
